{"js": "// 1. Remove the stray \"_GoBack\" bookmark that lived alone in the first,\n//    otherwise-empty paragraph of the document (it becomes a plain empty\n//    paragraph once the bookmark is gone).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. The \"For more information\" section used to end with:\n//      \"To stay informed about IBM training, see the following sites:\"\n//      \"IBM Training News: http://bit.ly/IBMTrainEN\"\n//      \"YouTube: youtube.com/IBMTraining\"\n//      \"Facebook: facebook.com/ibmtraining\"\n//      \"Twitter: twitter.com/websphere_edu\"\n//    All five paragraphs collapse into a single, empty paragraph that\n//    keeps the \"Abstract body text\" style / keepNext / keepLines\n//    formatting, plus the left indent (720 twips = 36pt) that the four\n//    \"label: link\" paragraphs had.\nconst body = context.document.body;\n\n// Locate the anchor paragraph by its distinctive text rather than by a\n// hard-coded index, so the script is resilient to other content in the\n// document.\nconst searchResults = body.search(\"To stay informed about IBM training\", {\n  matchCase: true,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\n\n// Walk forward and collect the next four paragraphs (IBM Training News /\n// YouTube / Facebook / Twitter), then delete them entirely.\nlet cursor = anchorParagraph;\nconst paragraphsToRemove = [];\nfor (let i = 0; i < 4; i++) {\n  cursor = cursor.getNext();\n  paragraphsToRemove.push(cursor);\n}\nparagraphsToRemove.forEach((paragraph) => paragraph.delete());\nawait context.sync();\n\n// Blank out the anchor paragraph's text and apply the indentation used by\n// the paragraphs that were just removed, leaving one empty paragraph\n// behind.\nanchorParagraph.clear();\nanchorParagraph.leftIndent = 36;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the stray \"_GoBack\" bookmark that lived alone in the first,\n#    otherwise-empty paragraph of the document (it becomes a plain empty\n#    paragraph once the bookmark is gone). \"_GoBack\" is a hidden bookmark,\n#    so it does not show up via Bookmarks.Count/iteration, but it can\n#    still be reached directly by name.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. The \"For more information\" section used to end with:\n#      \"To stay informed about IBM training, see the following sites:\"\n#      \"IBM Training News: http://bit.ly/IBMTrainEN\"\n#      \"YouTube: youtube.com/IBMTraining\"\n#      \"Facebook: facebook.com/ibmtraining\"\n#      \"Twitter: twitter.com/websphere_edu\"\n#    All five paragraphs collapse into a single, empty paragraph that\n#    keeps the \"Abstract body text\" style / keepNext / keepLines\n#    formatting, plus the left indent (720 twips = 36pt) that the four\n#    \"label: link\" paragraphs had.\n\n# Locate the anchor paragraph by its distinctive text rather than a\n# hard-coded index, so the script is resilient to other content in the\n# document.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"To stay informed about IBM training\")\nif (-not $found) {\n    throw \"Could not find the 'To stay informed about IBM training' paragraph\"\n}\n$anchorPos = $searchRange.Start\n\n$paragraphCount = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $paragraphCount; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Start -le $anchorPos -and $anchorPos -lt $candidate.Range.End) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not resolve the anchor paragraph index\"\n}\n\n$anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n$lastParagraphToDelete = $d.Paragraphs.Item($anchorIndex + 4)\n\n# Delete the four trailing \"label: link\" paragraphs (and their paragraph\n# marks) in one shot, from the end of the anchor paragraph through the end\n# of the fourth paragraph after it.\n$deleteRange = $d.Range($anchorParagraph.Range.End, $lastParagraphToDelete.Range.End)\n$deleteRange.Delete()\n\n# Re-fetch the (now final) anchor paragraph and blank out its own text\n# (but not its paragraph mark), then apply the indentation used by the\n# paragraphs that were just removed.\n$anchorParagraph = $d.Paragraphs.Item($anchorIndex)\n$textRange = $anchorParagraph.Range\n$textRange.End = $textRange.End - 1\n$textRange.Delete()\n$anchorParagraph.LeftIndent = 36\n"}
